$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing existing rows 104-127 down to 105-128.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the latest weekly price entry.
$ws.Cells.Item(104, 1).Value2  = 10
$ws.Cells.Item(104, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(104, 3).Value2  = "La Araucanía"
$ws.Cells.Item(104, 4).Value2  = 45173
$ws.Cells.Item(104, 5).Value2  = 9
$ws.Cells.Item(104, 6).Value2  = "Fruta"
$ws.Cells.Item(104, 7).Value2  = 100108
$ws.Cells.Item(104, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(104, 9).Value2  = 100108007
$ws.Cells.Item(104, 10).Value2 = "Coco"
$ws.Cells.Item(104, 11).Value2 = "Sin especificar"
$ws.Cells.Item(104, 12).Value2 = "Primera"
$ws.Cells.Item(104, 13).Value2 = 70
$ws.Cells.Item(104, 14).Value2 = 36000
$ws.Cells.Item(104, 15).Value2 = 36000
$ws.Cells.Item(104, 16).Value2 = 36000
$ws.Cells.Item(104, 17).Value2 = "$/malla 20 unidades"
$ws.Cells.Item(104, 18).Value2 = "Perú"
$ws.Cells.Item(104, 19).Value2 = 1800
$ws.Cells.Item(104, 20).Value2 = 20
